$d = $word.ActiveDocument

$q1 = [char]0x201C
$q2 = [char]0x201D

$find1 = $q1 + "B1x1_T-20.2_E35.00_FC" + $q2
$repl1 = $q1 + "B1x1_T-20_E35.0_FC" + $q2

$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false,
                         $true, 1, $false, $repl1, 2)

$find2 = $q1 + "B2x2_T-10.0_E180.0_FR"
$repl2 = $q1 + "B2x2_T-10_E180.0_FR"

$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false,
                         $true, 1, $false, $repl2, 2)
